$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.837.78'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').Value = '3.534.28'
$ws.Range('E3').Value = '  +1.07%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Formula = "'603.71"
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('E6').Value = '  +5.81%  '
$ws.Range('E7').Value = '  +0.61%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Formula = "'0.204"
$ws.Range('E9').Value = '  -2.69%  '
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('D11').Formula = "'53.73"
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('D12').Formula = "'0.0000303"
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('D13').Formula = "'9.53"
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('D14').Value = '4.090.63'
$ws.Range('E14').Value = '  +0.68%  '
$ws.Range('D15').Formula = "'604.53"
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').Value = '69.981.69'
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('D17').Formula = "'19.13"
$ws.Range('E17').Value = '  +0.79%  '
$ws.Range('D18').Formula = "'12.77"
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').Value = '3.511.42'
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('E20').Value = '  +1.00%  '
$ws.Range('D21').Formula = "'0.993"
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').Formula = "'18.22"
$ws.Range('E22').Value = '  +3.90%  '
$ws.Range('D23').Formula = "'5.31"
$ws.Range('E23').Value = '  +5.92%  '
$ws.Range('D24').Formula = "'102.21"
$ws.Range('E24').Value = '  -2.37%  '
$ws.Range('D25').Formula = "'4.62"
$ws.Range('E25').Value = '  -0.61%  '
$ws.Range('E26').Value = '  +4.45%  '
$ws.Range('E27').Value = '  +0.48%  '
$ws.Range('D28').Formula = "'9.66"
$ws.Range('E28').Value = '  -2.15%  '
$ws.Range('D29').Formula = "'33.52"
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').Formula = "'4.35"
$ws.Range('E30').Value = '  +16.67%  '
$ws.Range('E31').Value = '  +1.94%  '
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('E34').Value = '  -0.24%  '
$ws.Range('D35').Value = '0.0₃0867'
$ws.Range('E35').Value = '  +13.38%  '
$ws.Range('D36').Value = '3.739.51'
$ws.Range('E36').Value = '  +5.49%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').Formula = "'3.04"
$ws.Range('E38').Value = '  -2.37%  '
$ws.Range('D39').Formula = "'3.62"
$ws.Range('E39').Value = '  +1.41%  '
$ws.Range('D40').Formula = "'0.393"
$ws.Range('E40').Value = '  -0.62%  '
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('D42').Formula = "'485.16"
$ws.Range('E42').Value = '  -7.50%  '
$ws.Range('D43').Formula = "'0.134"
$ws.Range('E43').Value = '  -4.63%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('E45').Value = '  -3.50%  '
$ws.Range('E46').Value = '  -1.88%  '
$ws.Range('D47').Formula = "'3.32"
$ws.Range('E47').Value = '  -0.26%  '
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('D49').Formula = "'8.58"
$ws.Range('E49').Value = '  -3.15%  '
$ws.Range('D50').Formula = "'0.000254"
$ws.Range('E50').Value = '  +6.42%  '
$ws.Range('D51').Formula = "'130.57"
$ws.Range('E51').Value = '  -1.47%  '
